$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The structural edits below change formula/value dependents, so force a full
# recalculation the next time the workbook is opened (mirrors calcPr's
# fullCalcOnLoad flag that Excel stamps after this kind of edit).
$wb.ForceFullCalculation = $true

# Remove the "Xi (distance between pin wall to catheter wall)" column (old column I)
# this shifts old J,K,L (Yi, Mandrel Material, Mandrel OD) left into I,J,K
$ws.Range("I1").EntireColumn.Delete()

# Remove the old row 5 ("D" row, now empty of real data) entirely
$ws.Range("A5").EntireRow.Delete()

# Update row 2 (was id 5/"A") -> id 9/"A" with new measurements
$ws.Range("A2").Value = 9
$ws.Range("C2").Value = 300
$ws.Range("D2").Value = 1.6
$ws.Range("E2").Value = 1.2
$ws.Range("F2").Value = "Soft, black"
$ws.Range("I2").Value = 0.5

# Update row 3 (was id 6/"B") -> id 10/"B" with new measurements
$ws.Range("A3").Value = 10
$ws.Range("C3").Value = 800
$ws.Range("D3").Value = 1.6
$ws.Range("F3").Value = "Soft, black"
$ws.Range("H3").Value = 0.5
$ws.Range("I3").Value = 0.5

# Update row 4 (was id 7/"C") -> id 11/"C" with new measurements
$ws.Range("A4").Value = 11
$ws.Range("D4").Value = 1.6
$ws.Range("F4").Value = "Braided, purple"
$ws.Range("H4").Value = 0.5
$ws.Range("I4").Value = 0.5
